$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Row 4 corresponds to "Electricity" activity.
# Set "Exploiting Gas boiler for Heating" (H4), "Exploiting Gas boiler for
# Hot Sanitary Water" (K4), and "Exploiting Gas Stove for Cooking" (P4) to 0
# for the Electricity row - electricity use from gas boiler put to zero.
$ws.Range("H4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0
